# Generate Report for Handoff
# - Status changes from "In Translation" to "Ready for handoff"
# - Latest HO Xliff Generate Date / Handoff Datetime timestamps bump forward
# - Status column widened on all three sheets

$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$zhcn     = $wb.Worksheets.Item("zh-cn")
$dede     = $wb.Worksheets.Item("de-de")

# --- Update status text and timestamps -------------------------------------

# Overview sheet: zh-cn / de-de status columns (E, F) and the
# "Latest HO Xliff Generate Date" column (G)
$overview.Range("E2").Value = "Ready for handoff"
$overview.Range("F2").Value = "Ready for handoff"
$overview.Range("G2").Value = "2016-09-04 02:42:53"

# zh-cn sheet: Status column (C) and Latest Handoff Datetime (H)
$zhcn.Range("C2").Value = "Ready for handoff"
$zhcn.Range("H2").Value = "2016-09-04 02:42:49"

# de-de sheet: Status column (C) only (its handoff datetime cell already
# shares the same timestamp string that was updated via the Overview sheet)
$dede.Range("C2").Value = "Ready for handoff"

# --- Widen the Status / language status columns ----------------------------
# Target OOXML column width ~17.216; this runtime quantizes widths to
# 1/6-character steps, so 16.3 (-> stored width 17.1667, the nearest
# reachable step) is used to get as close as possible to the target.

$overview.Columns.Item(5).ColumnWidth = 16.3
$overview.Columns.Item(6).ColumnWidth = 16.3
$zhcn.Columns.Item(3).ColumnWidth = 16.3
$dede.Columns.Item(3).ColumnWidth = 16.3
